$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (header already in row 1).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.
$data = @(
    @("FAPs","Rspo3","Fzd8","ECs",3,1,2.507621333333333,7.522864,0.9939780200440224,0.9939780200440224,2,0.6666666666666666,2.558821666666667,7.676465,0.2156728774407755,0.2156728774407755,6.416555799528889,57.74900219576001,0.2143740996957791,0.2143740996957791),
    @("FAPs","Rspo3","Fzd8","FAPs",3,1,2.507621333333333,7.522864,0.9939780200440224,0.9939780200440224,3,1,6.453984666666667,19.361954,0.5439806384912759,0.5439806384912759,16.18414963513955,145.657346716256,0.5407047979898415,0.5407047979898415),
    @("FAPs","Rspo3","Fzd8","sCs",3,1,2.507621333333333,7.522864,0.9939780200440224,0.9939780200440224,3,1,2.851558333333334,8.554675000000001,0.2403464840679487,0.2403464840679487,7.150628509911112,64.35565658920001,0.2388991223584018,0.2388991223584019),
    @("sCs","Rspo3","Fzd8","ECs",1,0.3333333333333333,0.01519233333333333,0.045577,0.006021979955977724,0.006021979955977723,2,0.6666666666666666,2.558821666666667,7.676465,0.2156728774407755,0.2156728774407755,0.03887447170055556,0.349870245305,0.00129877774499639,0.00129877774499639),
    @("sCs","Rspo3","Fzd8","FAPs",1,0.3333333333333333,0.01519233333333333,0.045577,0.006021979955977724,0.006021979955977723,3,1,6.453984666666667,19.361954,0.5439806384912759,0.5439806384912759,0.09805108638422223,0.882459777458,0.003275840501434428,0.003275840501434427),
    @("sCs","Rspo3","Fzd8","sCs",1,0.3333333333333333,0.01519233333333333,0.045577,0.006021979955977724,0.006021979955977723,3,1,2.851558333333334,8.554675000000001,0.2403464840679487,0.2403464840679487,0.04332182471944445,0.3898964224750001,0.001447361709546907,0.001447361709546907)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
